$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing sheet so we get a second tire sheet ("2x multi-axle" support)
$ws1.Copy($null, $ws1)

$wsOrig = $wb.Worksheets.Item(1)
$wsCopy = $wb.Worksheets.Item(2)

# Rename sheets to the new naming convention
$wsOrig.Name = "Tire2x_270_70R22"
$wsCopy.Name = "Tire2x_430_50R38"

# The first (original) sheet becomes the "270_70R22" tire and gets a fixed numeric
# value for the half-width offset instead of the old formula, and its title label
# text is updated to match the new sheet name.
$wsOrig.Range("H7").Value = 0.4572
$wsOrig.Range("H2").Value = "Tire"

# The new sheet keeps the old formula/values, but its title label should read the
# new sheet's own name. Set this before the original sheet's H3 label so that the
# shared-string table keeps the same insertion order as Excel produced.
$wsCopy.Range("H3").Value = "Tire2x_430_50R38"
$wsOrig.Range("H3").Value = "Tire2x_270_70R22"

$wsOrig.Range("C25").Select()

# The new sheet becomes the active/selected sheet.
$wsCopy.Select()
$wsCopy.Range("J16").Select()
